$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case connector words (de/del/el/la/las/los/y) in state/municipality names, plus two special-case fixes ---
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B16").Value = "Playas De Rosarito"
$ws.Range("B40").Value = "Amatenango De La Frontera"
$ws.Range("B41").Value = "Amatenango Del Valle"
$ws.Range("B44").Value = "Bejucal De Ocampo"
$ws.Range("B46").Value = "Benemérito De Las Américas"
$ws.Range("B54").Value = "Chiapa De Corzo"
$ws.Range("B59").Value = "Comitán De Domínguez"
$ws.Range("B86").Value = "Mazapa De Madero"
$ws.Range("B89").Value = "Montecristo De Guerero"
$ws.Range("B93").Value = "Ocozocoautla De Espinosa"
$ws.Range("B104").Value = "Salto De Agua"
$ws.Range("B105").Value = "San Cristóbal De Las Casas"
$ws.Range("B150").Value = "Coyame Del Sotol"
$ws.Range("B160").Value = "Guadalupe Y Calvo"
$ws.Range("B163").Value = "Hidalgo Del Parral"
$ws.Range("B185").Value = "San Francisco De Borja"
$ws.Range("B186").Value = "San Francisco De Conchos"
$ws.Range("B187").Value = "San Francisco Del Oro"
$ws.Range("B195").Value = "Valle De Zaragoza"
$ws.Range("B216").Value = "San Juan De Sabinas"
$ws.Range("B231").Value = "Villa De Álvarez"
$ws.Range("A233").Value = "Ciudad De México"
$ws.Range("B237").Value = "Cuajimalpa De Morelos"
$ws.Range("B252").Value = "Coneto De Comonfort"
$ws.Range("B266").Value = "Nombre De Dios"
$ws.Range("B270").Value = "Pánuco De Coronado"
$ws.Range("B277").Value = "San Juan De Guadalupe"
$ws.Range("B278").Value = "San Juan Del Río"
$ws.Range("B279").Value = "San Luis Del Cordero"
$ws.Range("B280").Value = "San Pedro Del Gallo"
$ws.Range("A290").Value = "Estado De México"
$ws.Range("B290").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B293").Value = "Almoloya De Alquisiras"
$ws.Range("B294").Value = "Almoloya De Juárez"
$ws.Range("B295").Value = "Almoloya Del Río"
$ws.Range("B302").Value = "Atizapán De Zaragoza"
$ws.Range("B310").Value = "Chapa De Mota"
$ws.Range("B316").Value = "Coacalco De Berriozábal"
$ws.Range("B323").Value = "Ecatepec De Morelos"
$ws.Range("B331").Value = "Ixtapan De La Sal"
$ws.Range("B332").Value = "Ixtapan Del Oro"
$ws.Range("B349").Value = "Naucalpan De Juárez"
$ws.Range("B363").Value = "San Antonio La Isla"
$ws.Range("B364").Value = "San Felipe Del Progreso"
$ws.Range("B365").Value = "San Martín De Las Pirámides"
$ws.Range("B367").Value = "San Simón De Guerero"
$ws.Range("B369").Value = "Soyaniquilpan De Juárez"
$ws.Range("B379").Value = "Tenango Del Aire"
$ws.Range("B380").Value = "Tenango Del Valle"
$ws.Range("B394").Value = "Tlalnepantla De Baz"
$ws.Range("B400").Value = "Valle De Bravo"
$ws.Range("B401").Value = "Valle De Chalco Solidaridad"
$ws.Range("B402").Value = "Villa De Allende"
$ws.Range("B403").Value = "Villa Del Carbón"
$ws.Range("B417").Value = "Apaseo El Alto"
$ws.Range("B418").Value = "Apaseo El Grande"
$ws.Range("B426").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B430").Value = "Jaral Del Progreso"
$ws.Range("B438").Value = "Purísima Del Rincón"
$ws.Range("B442").Value = "San Diego De La Unión"
$ws.Range("B444").Value = "San Francisco Del Rincón"
$ws.Range("B446").Value = "San Luis De La Paz"
$ws.Range("B448").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B450").Value = "Silao De La Victoria"
$ws.Range("B455").Value = "Valle De Santiago"
$ws.Range("B461").Value = "Acapulco De Juárez"
$ws.Range("B464").Value = "Ajuchitlán Del Progreso"
$ws.Range("B465").Value = "Alcozauca De Guerero"
$ws.Range("B469").Value = "Atenango Del Río"
$ws.Range("B470").Value = "Atlamajalcingo Del Monte"
$ws.Range("B472").Value = "Atoyac De Álvarez"
$ws.Range("B473").Value = "Ayutla De Los Libres"
$ws.Range("B476").Value = "Buenavista De Cuéllar"
$ws.Range("B477").Value = "Chilapa De Álvarez"
$ws.Range("B478").Value = "Chilpancingo De Los Bravo"
$ws.Range("B479").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B484").Value = "Coyuca De Benítez"
$ws.Range("B485").Value = "Coyuca De Catalán"
$ws.Range("B489").Value = "Cuetzala Del Progreso"
$ws.Range("B490").Value = "Cutzamala De Pinzón"
$ws.Range("B496").Value = "Huitzuco De Los Figueroa"
$ws.Range("B497").Value = "Iguala De La Independencia"
$ws.Range("B499").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B500").Value = "Zihuatanejo De Azueta"
$ws.Range("B502").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B505").Value = "Mártir De Cuilapan"
$ws.Range("B518").Value = "Taxco De Alarcón"
$ws.Range("B520").Value = "Técpan De Galeana"
$ws.Range("B522").Value = "Tepecoacuilco De Trujano"
$ws.Range("B524").Value = "Tixtla De Guerero"
$ws.Range("B528").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B529").Value = "Tlapa De Comonfort"
$ws.Range("B541").Value = "Agua Blanca De Iturbide"
$ws.Range("B548").Value = "Atotonilco De Tula"
$ws.Range("B549").Value = "Atotonilco El Grande"
$ws.Range("B555").Value = "Cuautepec De Hinojosa"
$ws.Range("B561").Value = "Huasca De Ocampo"
$ws.Range("B565").Value = "Huejutla De Reyes"
$ws.Range("B568").Value = "Jacala De Ledezma"
$ws.Range("B575").Value = "Mineral De La Reforma"
$ws.Range("B576").Value = "Mineral Del Chico"
$ws.Range("B577").Value = "Mineral Del Monte"
$ws.Range("B578").Value = "Mixquiahuala De Juárez"
$ws.Range("B579").Value = "Molango De Escamilla"
$ws.Range("B581").Value = "Nopala De Villagrán"
$ws.Range("B582").Value = "Omitlán De Juárez"
$ws.Range("B583").Value = "Pachuca De Soto"
$ws.Range("B586").Value = "Progreso De Obregón"
$ws.Range("B592").Value = "Santiago De Anaya"
$ws.Range("B593").Value = "Santiago Tulantepec De Lugo Guerero"
$ws.Range("B597").Value = "Tenango De Doria"
$ws.Range("B599").Value = "Tepehuacán De Guerero"
$ws.Range("B600").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B603").Value = "Tezontepec De Aldama"
$ws.Range("B612").Value = "Tula De Allende"
$ws.Range("B613").Value = "Tulancingo De Bravo"
$ws.Range("B614").Value = "Villa De Tezontepec"
$ws.Range("B618").Value = "Zacualtipán De Ángeles"
$ws.Range("B619").Value = "Zapotlán De Juárez"
$ws.Range("B624").Value = "Acatlán De Juárez"
$ws.Range("B625").Value = "Ahualulco De Mercado"
$ws.Range("B630").Value = "Atemajac De Brizuela"
$ws.Range("B633").Value = "Atotonilco El Alto"
$ws.Range("B635").Value = "Autlán De Navarro"
$ws.Range("B641").Value = "Cañadas De Obregón"
$ws.Range("B648").Value = "Concepción De Buenos Aires"
$ws.Range("B649").Value = "Cuautitlán De García Barragán"
$ws.Range("B658").Value = "Encarnación De Díaz"
$ws.Range("B665").Value = "Huejuquilla El Alto"
$ws.Range("B666").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B667").Value = "Ixtlahuacán Del Río"
$ws.Range("B671").Value = "Jilotlán De Los Dolores"
$ws.Range("B677").Value = "La Manzanilla De La Paz"
$ws.Range("B678").Value = "Lagos De Moreno"
$ws.Range("B686").Value = "Ojuelos De Jalisco"
$ws.Range("B691").Value = "San Cristóbal De La Barranca"
$ws.Range("B692").Value = "San Diego De Alejandría"
$ws.Range("B694").Value = "San Juan De Los Lagos"
$ws.Range("B695").Value = "San Juanito De Escobedo"
$ws.Range("B698").Value = "San Martín De Bolaños"
$ws.Range("B700").Value = "San Miguel El Alto"
$ws.Range("B701").Value = "San Sebastián Del Oeste"
$ws.Range("B702").Value = "Santa María De Los Ángeles"
$ws.Range("B703").Value = "Santa María Del Oro"
$ws.Range("B706").Value = "Talpa De Allende"
$ws.Range("B707").Value = "Tamazula De Gordiano"
$ws.Range("B710").Value = "Techaluta De Montenegro"
$ws.Range("B714").Value = "Teocuitatlán De Corona"
$ws.Range("B715").Value = "Tepatitlán De Morelos"
$ws.Range("B718").Value = "Tizapán El Alto"
$ws.Range("B719").Value = "Tlajomulco De Zúñiga"
$ws.Range("B731").Value = "Unión De San Antonio"
$ws.Range("B732").Value = "Unión De Tula"
$ws.Range("B733").Value = "Valle De Guadalupe"
$ws.Range("B734").Value = "Valle De Juárez"
$ws.Range("B739").Value = "Yahualica De González Gallo"
$ws.Range("B740").Value = "Zacoalco De Torres"
$ws.Range("B743").Value = "Zapotitlán De Vadillo"
$ws.Range("B744").Value = "Zapotlán Del Rey"
$ws.Range("B745").Value = "Zapotlán El Grande"
$ws.Range("B771").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B773").Value = "Cojumatlán De Régules"
$ws.Range("B840").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B866").Value = "Coatlán Del Río"
$ws.Range("B874").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B878").Value = "Puente De Ixtla"
$ws.Range("B884").Value = "Tetela Del Volcán"
$ws.Range("B886").Value = "Tlaltizapán De Zapata"
$ws.Range("B894").Value = "Zacualpan De Amilpas"
$ws.Range("B898").Value = "Amatlán De Cañas"
$ws.Range("B899").Value = "Bahía De Banderas"
$ws.Range("B903").Value = "Ixtlán Del Río"
$ws.Range("B910").Value = "Santa María Del Oro"
$ws.Range("B930").Value = "Mier Y Noriega"
$ws.Range("B935").Value = "San Nicolás De Los Garza"
$ws.Range("B940").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B948").Value = "Ayoquezco De Aldama"
$ws.Range("B952").Value = "Capulálpam De Méndez"
$ws.Range("B954").Value = "Chalcatongo De Hidalgo"
$ws.Range("B955").Value = "Ciénega De Zimatlán"
$ws.Range("B958").Value = "Coicoyán De Las Flores"
$ws.Range("B961").Value = "Constancia Del Rosario"
$ws.Range("B964").Value = "Cuilápam De Guerero"
$ws.Range("B965").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B966").Value = "El Barrio De La Soledad"
$ws.Range("B968").Value = "Eloxochitlán De Flores Magón"
$ws.Range("B969").Value = "Fresnillo De Trujano"
$ws.Range("B970").Value = "Guadalupe De Ramírez"
$ws.Range("B972").Value = "Guelatao De Juárez"
$ws.Range("B973").Value = "Guevea De Humboldt"
$ws.Range("B974").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B975").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B976").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B978").Value = "Huautla De Jiménez"
$ws.Range("B980").Value = "Ixtlán De Juárez"
$ws.Range("B981").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B995").Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Range("B997").Value = "Mariscala De Juárez"
$ws.Range("B998").Value = "Mártires De Tacubaya"
$ws.Range("B1000").Value = "Mazatlán Villa De Flores"
$ws.Range("B1002").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B1003").Value = "Mixistlán De La Reforma"
$ws.Range("B1007").Value = "Nejapa De Madero"
$ws.Range("B1009").Value = "Oaxaca De Juárez"
$ws.Range("B1010").Value = "Ocotlán De Morelos"
$ws.Range("B1011").Value = "Pinotepa De Don Luis"
$ws.Range("B1013").Value = "Putla Villa De Guerero"
$ws.Range("B1014").Value = "Reforma De Pineda"
$ws.Range("B1016").Value = "Rojas De Cuauhtémoc"
$ws.Range("B1021").Value = "San Agustín De Las Juntas"
$ws.Range("B1040").Value = "San Antonino El Alto"
$ws.Range("B1043").Value = "San Antonio De La Cal"
$ws.Range("B1050").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B1066").Value = "San Dionisio Del Mar"
$ws.Range("B1070").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B1077").Value = "San Francisco Del Mar"
$ws.Range("B1100").Value = "San José Del Peñasco"
$ws.Range("B1101").Value = "San José Del Progreso"
$ws.Range("B1113").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B1127").Value = "San Juan De Los Cués"
$ws.Range("B1128").Value = "San Juan Del Estado"
$ws.Range("B1129").Value = "San Juan Del Río"
$ws.Range("B1170").Value = "San Martín De Los Cansecos"
$ws.Range("B1178").Value = "San Mateo Del Mar"
$ws.Range("B1192").Value = "San Miguel Del Puerto"
$ws.Range("B1193").Value = "San Miguel Del Río"
$ws.Range("B1195").Value = "San Miguel El Grande"
$ws.Range("B1219").Value = "San Pablo Villa De Mitla"
$ws.Range("B1227").Value = "San Pedro El Alto"
$ws.Range("B1248").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1249").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B1250").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B1268").Value = "Santa Ana Del Valle"
$ws.Range("B1285").Value = "Santa Cruz De Bravo"
$ws.Range("B1290").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B1296").Value = "Santa Inés De Zaragoza"
$ws.Range("B1297").Value = "Santa Inés Del Monte"
$ws.Range("B1299").Value = "Santa Lucía Del Camino"
$ws.Range("B1313").Value = "Santa María Del Rosario"
$ws.Range("B1314").Value = "Santa María Del Tule"
$ws.Range("B1322").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1361").Value = "Santiago Del Río"
$ws.Range("B1399").Value = "Santo Domingo De Morelos"
$ws.Range("B1424").Value = "Sitio De Xitlapehua"
$ws.Range("B1426").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1427").Value = "Tanetze De Zaragoza"
$ws.Range("B1429").Value = "Tataltepec De Valdés"
$ws.Range("B1430").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B1431").Value = "Teotitlán De Flores Magón"
$ws.Range("B1432").Value = "Teotitlán Del Valle"
$ws.Range("B1434").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B1435").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1436").Value = "Tlacolula De Matamoros"
$ws.Range("B1438").Value = "Tlalixtac De Cabrera"
$ws.Range("B1439").Value = "Totontepec Villa De Morelos"
$ws.Range("B1443").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1444").Value = "Villa De Etla"
$ws.Range("B1445").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1446").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1447").Value = "Villa De Zaachila"
$ws.Range("B1450").Value = "Villa Sola De Vega"
$ws.Range("B1451").Value = "Villa Talea De Castro"
$ws.Range("B1454").Value = "Yutanduchi De Guerero"
$ws.Range("B1455").Value = "Zapotitlán Del Río"
$ws.Range("B1458").Value = "Zimatlán De Álvarez"
$ws.Range("B1485").Value = "Ayotoxco De Guerero"
$ws.Range("B1489").Value = "Chalchicomula De Sesma"
$ws.Range("B1499").Value = "Chila De La Sal"
$ws.Range("B1510").Value = "Cuapiaxtla De Madero"
$ws.Range("B1514").Value = "Cuayuca De Andrade"
$ws.Range("B1515").Value = "Cuetzalan Del Progreso"
$ws.Range("B1531").Value = "Huehuetlán El Chico"
$ws.Range("B1532").Value = "Huehuetlán El Grande"
$ws.Range("B1537").Value = "Huitzilan De Serdán"
$ws.Range("B1539").Value = "Ixcamilpa De Guerero"
$ws.Range("B1542").Value = "Izúcar De Matamoros"
$ws.Range("B1552").Value = "Los Reyes De Juárez"
$ws.Range("B1553").Value = "Mazapiltepec De Juárez"
$ws.Range("B1566").Value = "Palmar De Bravo"
$ws.Range("B1576").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1593").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1597").Value = "San Salvador El Seco"
$ws.Range("B1598").Value = "San Salvador El Verde"
$ws.Range("B1606").Value = "Tecali De Herrera"
$ws.Range("B1614").Value = "Tepanco De López"
$ws.Range("B1615").Value = "Tepango De Rodríguez"
$ws.Range("B1616").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1622").Value = "Tepexi De Rodríguez"
$ws.Range("B1624").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1625").Value = "Tetela De Ocampo"
$ws.Range("B1626").Value = "Teteles De Avila Castillo"
$ws.Range("B1631").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1643").Value = "Totoltepec De Guerero"
$ws.Range("B1645").Value = "Tuzamapan De Galeana"
$ws.Range("B1649").Value = "Xayacatlán De Bravo"
$ws.Range("B1655").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1663").Value = "Zapotitlán De Méndez"
$ws.Range("B1672").Value = "Amealco De Bonfil"
$ws.Range("B1674").Value = "Cadereyta De Montes"
$ws.Range("B1680").Value = "Jalpan De Serra"
$ws.Range("B1681").Value = "Landa De Matamoros"
$ws.Range("B1684").Value = "Pinal De Amoles"
$ws.Range("B1687").Value = "San Juan Del Río"
$ws.Range("B1699").Value = "Armadillo De Los Infante"
$ws.Range("B1700").Value = "Axtla De Terrazas"
$ws.Range("B1706").Value = "Ciudad Del Maíz"
$ws.Range("B1715").Value = "Mexquitic De Carmona"
$ws.Range("B1721").Value = "San Ciro De Acosta"
$ws.Range("B1724").Value = "Santa María Del Río"
$ws.Range("B1726").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1734").Value = "Villa De Arista"
$ws.Range("B1735").Value = "Villa De Arriaga"
$ws.Range("B1736").Value = "Villa De Guadalupe"
$ws.Range("B1737").Value = "Villa De Ramos"
$ws.Range("B1738").Value = "Villa De Reyes"
$ws.Range("B1799").Value = "Nacozari De García"
$ws.Range("B1811").Value = "San Pedro De La Cueva"
$ws.Range("B1828").Value = "Jalpa De Méndez"
$ws.Range("B1862").Value = "Soto La Marina"
$ws.Range("B1869").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1871").Value = "Amaxac De Guerero"
$ws.Range("B1872").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1878").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1886").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1890").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1891").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1892").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1895").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1898").Value = "San Pablo Del Monte"
$ws.Range("B1899").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B1906").Value = "Tepetitla De Lardizábal"
$ws.Range("B1909").Value = "Tetla De La Solidaridad"
$ws.Range("B1920").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1930").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1934").Value = "Amatlán De Los Reyes"
$ws.Range("B1946").Value = "Boca Del Río"
$ws.Range("B1948").Value = "Camarón De Tejeda"
$ws.Range("B1952").Value = "Castillo De Teayo"
$ws.Range("B1954").Value = "Cazones De Herrera"
$ws.Range("B1962").Value = "Chinampa De Gorostiza"
$ws.Range("B1975").Value = "Cosamaloapan De Carpio"
$ws.Range("B1976").Value = "Cosautlán De Carvajal"
$ws.Range("B1992").Value = "Hueyapan De Ocampo"
$ws.Range("B1993").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1994").Value = "Ignacio De La Llave"
$ws.Range("B1998").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1999").Value = "Ixhuatlán De Madero"
$ws.Range("B2000").Value = "Ixhuatlán Del Café"
$ws.Range("B2001").Value = "Ixhuatlán Del Sureste"
$ws.Range("B2010").Value = "Juchique De Ferrer"
$ws.Range("B2014").Value = "Las Vigas De Ramírez"
$ws.Range("B2015").Value = "Lerdo De Tejada"
$ws.Range("B2020").Value = "Martínez De La Torre"
$ws.Range("B2023").Value = "Medellín De Bravo"
$ws.Range("B2028").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B2039").Value = "Ozuluama De Mascareñas"
$ws.Range("B2043").Value = "Paso De Ovejas"
$ws.Range("B2044").Value = "Paso Del Macho"
$ws.Range("B2048").Value = "Poza Rica De Hidalgo"
$ws.Range("B2059").Value = "Sayula De Alemán"
$ws.Range("B2062").Value = "Soledad De Doblado"
$ws.Range("B2068").Value = "Tatahuicapan De Juárez"
$ws.Range("B2103").Value = "Vega De Alatorre"
$ws.Range("B2115").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B2116").Value = "Zozocolco De Hidalgo"
$ws.Range("B2170").Value = "Tekal De Venegas"
$ws.Range("B2197").Value = "Cañitas De Felipe Pescador"
$ws.Range("B2199").Value = "Concepción Del Oro"
$ws.Range("B2201").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B2212").Value = "Jiménez Del Teul"
$ws.Range("B2218").Value = "Mezquital Del Oro"
$ws.Range("B2223").Value = "Moyahua De Estrada"
$ws.Range("B2224").Value = "Nochistlán De Mejía"
$ws.Range("B2225").Value = "Noria De Ángeles"
$ws.Range("B2236").Value = "Teúl De González Ortega"
$ws.Range("B2237").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B2239").Value = "Trinidad García De La Cadena"
$ws.Range("B2242").Value = "Villa De Cos"
$ws.Range("A414").Value = "Guanajuato"
$ws.Range("B931").Value = "Montemorelos"

# --- Floating point literal fixes for D column (Excel recalculation last-bit rounding) ---
$ws.Range("D309").Value = 0.000958056835312848
$ws.Range("D480").Value = 0.0009017005508826804
$ws.Range("D491").Value = 0.0009228341575439932
$ws.Range("D502").Value = 0.000958056835312848
$ws.Range("D636").Value = 0.000943967764205306
$ws.Range("D651").Value = 0.0009017005508826804
$ws.Range("D746").Value = 0.000943967764205306
$ws.Range("D755").Value = 0.000958056835312848
$ws.Range("D763").Value = 0.0009510122997590768
$ws.Range("D787").Value = 0.0009510122997590768
$ws.Range("D795").Value = 0.0009087450864364512
$ws.Range("D815").Value = 0.0009510122997590768
$ws.Range("D1056").Value = 0.0009721459064203896
$ws.Range("D1067").Value = 0.0009721459064203896
$ws.Range("D1365").Value = 0.0009721459064203896
$ws.Range("D1462").Value = 0.000958056835312848
$ws.Range("D1587").Value = 0.000943967764205306
$ws.Range("D1616").Value = 0.0009228341575439932
$ws.Range("D1672").Value = 0.0009721459064203896
$ws.Range("D1722").Value = 0.0009721459064203896
$ws.Range("D1760").Value = 0.000943967764205306
$ws.Range("D1911").Value = 0.000943967764205306

# --- Remove trailing footer/metadata rows (2250-2255), shrinking used range to A1:D2249 ---
$ws.Range("A2250:D2255").ClearContents()

